$d = $word.ActiveDocument

function Format-TagRun($rng) {
    # Re-style a narrow sub-range (the literal "<exp>" / "</exp>" marker text)
    # as small gray Courier New, matching the editorial-markup convention
    # used elsewhere in this document.
    $rng.Find.ClearFormatting()
    $rng.Find.Replacement.ClearFormatting()
    $rng.Find.Replacement.Font.Name = "Courier New"
    $rng.Find.Replacement.Font.NameFarEast = "Courier New"
    $rng.Find.Replacement.Font.NameBi = "Courier New"
    $rng.Find.Replacement.Font.Color = 11119017
    $rng.Find.Replacement.Font.Size = 7
    $txt = $rng.Text
    $rng.Find.Execute($txt, $false, $false, $false, $false, $false, $true, 1, $true, $txt, 2) | Out-Null
}

function Split-ExpTags($matchText, [int[]]$openStarts, [int[]]$openLens, [int[]]$closeStarts, [int[]]$closeLens) {
    # not used; kept simple inline below
}

# ---------------------------------------------------------------------
# Hunk 1: " en baigna<exp>n</exp>t"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Replacement.ClearFormatting()
$found = $rng.Find.Execute(" en baigna<exp>n</exp>t", $false)
if (-not $found) { throw "hunk1 text not found" }
$start = $rng.Start
Format-TagRun ($d.Range($start+10, $start+15))   # "<exp>"
Format-TagRun ($d.Range($start+16, $start+22))   # "</exp>"

# ---------------------------------------------------------------------
# Hunk 3: "premierem<exp>ent</exp> à la "
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Replacement.ClearFormatting()
$found = $rng.Find.Execute("premierem<exp>ent</exp> à la ", $false)
if (-not $found) { throw "hunk3 text not found" }
$start = $rng.Start
Format-TagRun ($d.Range($start+9, $start+14))    # "<exp>"
Format-TagRun ($d.Range($start+17, $start+23))   # "</exp>"

# ---------------------------------------------------------------------
# Hunk 4: ", pas trop allumée pour donner feu lent au coma<exp>n</exp>cem<exp>ent</exp>,"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Replacement.ClearFormatting()
$found = $rng.Find.Execute(", pas trop allumée pour donner feu lent au coma<exp>n</exp>cem<exp>ent</exp>,", $false)
if (-not $found) { throw "hunk4 text not found" }
$start = $rng.Start
Format-TagRun ($d.Range($start+47, $start+52))   # "<exp>"
Format-TagRun ($d.Range($start+53, $start+59))   # "</exp>"
Format-TagRun ($d.Range($start+62, $start+67))   # "<exp>"
Format-TagRun ($d.Range($start+70, $start+76))   # "</exp>"

# ---------------------------------------------------------------------
# Hunk 5: "à bon foeu, mesmem<exp>ent</exp> ceulx de"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Replacement.ClearFormatting()
$found = $rng.Find.Execute("à bon foeu, mesmem<exp>ent</exp> ceulx de", $false)
if (-not $found) { throw "hunk5 text not found" }
$start = $rng.Start
Format-TagRun ($d.Range($start+18, $start+23))   # "<exp>"
Format-TagRun ($d.Range($start+26, $start+32))   # "</exp>"

# ---------------------------------------------------------------------
# Hunk 6: " ne se fendro<exp>n</exp>t"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Replacement.ClearFormatting()
$found = $rng.Find.Execute(" ne se fendro<exp>n</exp>t", $false)
if (-not $found) { throw "hunk6 text not found" }
$start = $rng.Start
Format-TagRun ($d.Range($start+13, $start+18))   # "<exp>"
Format-TagRun ($d.Range($start+19, $start+25))   # "</exp>"

# ---------------------------------------------------------------------
# Hunk 7: "Recuits doulcem<exp>ent</exp>"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Replacement.ClearFormatting()
$found = $rng.Find.Execute("Recuits doulcem<exp>ent</exp>", $false)
if (-not $found) { throw "hunk7 text not found" }
$start = $rng.Start
Format-TagRun ($d.Range($start+15, $start+20))   # "<exp>"
Format-TagRun ($d.Range($start+23, $start+29))   # "</exp>"

# ---------------------------------------------------------------------
# Hunk 8: "couches, mesmem<exp>ent</exp> au"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Replacement.ClearFormatting()
$found = $rng.Find.Execute("couches, mesmem<exp>ent</exp> au", $false)
if (-not $found) { throw "hunk8 text not found" }
$start = $rng.Start
Format-TagRun ($d.Range($start+15, $start+20))   # "<exp>"
Format-TagRun ($d.Range($start+23, $start+29))   # "</exp>"

# ---------------------------------------------------------------------
# Hunk 2: " sont plustost recuits, prens " -> " sont plustost recuits. Prens "
#   the "." and "P" land in runs with no explicit color (rtl-only rPr)
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Replacement.ClearFormatting()
$found = $rng.Find.Execute(" sont plustost recuits, prens ", $false)
if (-not $found) { throw "hunk2 text not found" }
$start = $rng.Start

# "," (single char at offset 22) -> "."
$commaRng = $d.Range($start+22, $start+23)
$commaRng.Find.ClearFormatting()
$commaRng.Find.Replacement.ClearFormatting()
$commaRng.Find.Replacement.Font.ColorIndex = 0
$commaRng.Find.Execute(",", $false, $false, $false, $false, $false, $true, 1, $true, ".", 2) | Out-Null

# "p" (single char at offset 24) -> "P"
$pRng = $d.Range($start+24, $start+25)
$pRng.Find.ClearFormatting()
$pRng.Find.Replacement.ClearFormatting()
$pRng.Find.Replacement.Font.ColorIndex = 0
$pRng.Find.Execute("p", $false, $false, $false, $false, $false, $true, 1, $true, "P", 2) | Out-Null
